# Manual Renter All Stocks
# - Adds "long straddle" / "Long Straddle Status" columns (E/F) next to the
#   existing "short straddle" column (D) on Sheet1.
# - Widens/adds the new columns B, C, E, F to match the target layout.
# - Moves the active selection to D12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (these create two new shared-string entries: "long
# straddle" and "Long Straddle Status", and extend row 1's used range).
$ws.Range("E1").Value = "long straddle"
$ws.Range("F1").Value = "Long Straddle Status"

# Column widths (B, C, E, F are newly custom-sized; A and D already had a
# custom width from the original workbook and are left untouched).
$ws.Columns.Item(2).ColumnWidth = 24
$ws.Columns.Item(3).ColumnWidth = 46.666666666666664
$ws.Columns.Item(5).ColumnWidth = 30.333333333333332
$ws.Columns.Item(6).ColumnWidth = 17

# Move / restore the selection to D12, matching the saved view state.
$ws.Range("D12").Select() | Out-Null
